# Fix: when all observations in an aggregating group are NA, sum(na.rm=T)
# previously returned 0; it should now return NA. Update the affected
# "Value" (column G) cells on the "Data" sheet from "0" to "NaN".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$rows = @(12, 19, 20, 23, 301, 302, 584, 585)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7 = "Value"
    if ($cell.Value2 -eq "0") {
        $cell.Value = "NaN"
    }
}
